$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Update existing "blue" cell to "blue -1"
$ws1.Range("A2").Value = "blue -1"

# Add new row with "orange "
$ws1.Range("A3").Value = "orange "

# Add a new sheet "Sheet2" after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"
$ws2.Range("A1").Value = "new sheet "

$ws1.Select()
$ws1.Range("A3").Select()
